# "Fertigstellung d. Content Analyse"
# Fill in the previously-empty "Content" description cells (column H for the
# "Ebene 2" rows, column L for the "Ebene 3" rows) for rows 172-206, and
# restore the window's scroll position / active-cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('H172').Value2 = 'Literatur, Antragsformular und Links zu High Performance Computing'
$ws.Range('H173').Value2 = 'Information, Kontaktdaten, personalisierte SW-Lösung'
$ws.Range('L174').Value2 = 'spezielle Information zu EVO und Panda, Links'
$ws.Range('H175').Value2 = 'Sophos, Installation, Aktualisierung und Konfiguration von Virenscanner'
$ws.Range('H176').Value2 = 'Information zu Virtualisierung, Produktübersicht, Anforderungen., Vorteile, Datenaustausch, Installation'
$ws.Range('H177').Value2 = 'Kontaktdaten, Schwarzes Brett, Zugang Webserver, Dokumentation'
$ws.Range('L178').Value2 = 'allg. Nutzungsbedingungen für private Homepages'
$ws.Range('L179').Value2 = 'FAQs'
$ws.Range('L180').Value2 = 'Antrag eigene Homepage, Zugangsmethoden'
$ws.Range('L181').Value2 = 'versch. Information zu persönlicher Homepage'
$ws.Range('L182').Value2 = 'Zugangsmethoden zur eigenen Homepage'
$ws.Range('L183').Value2 = 'Information, Antrag auf Zugang, Anleitung Samba / Linux, Zugriffssteuerung'
$ws.Range('L184').Value2 = 'Hilfestellung: Dateinamen Konvention'
$ws.Range('H185').Value2 = 'Anleitung automatisierte Windows-Installation'
$ws.Range('H186').Value2 = 'externer Link (RZ-Windows-Update-Service)'
$ws.Range('H187').Value2 = 'ReX-Linux d. Uni R, Allgemeine Information'
$ws.Range('L188').Value2 = 'Installationsanleitung, Partitionierung'
$ws.Range('L189').Value2 = 'Administrationsanleitung'
$ws.Range('H190').Value2 = 'leere Seite'
$ws.Range('L191').Value2 = 'Download TSM-Client'
$ws.Range('H193').Value2 = 'Kontaktdaten'
$ws.Range('H194').Value2 = 'allgemeine Information zu Ausbildungsstellen und Kontaktdaten'
$ws.Range('L195').Value2 = 'Informationen Azubi-Camp'
$ws.Range('L196').Value2 = 'Aufgaben d. Azubis'
$ws.Range('H197').Value2 = 'Auflistung aller Mitarbeiter und Kontaktdaten'
$ws.Range('H198').Value2 = 'Anfahrtgsbeschreibung, Gebäudepläne'
$ws.Range('H199').Value2 = 'aktuelle Stellenausschribung, Kontaktdaten, Arbeits-/Entgeltbedingungen'
$ws.Range('H200').Value2 = 'Anmeldung zum Newsletter (RZettel)'
$ws.Range('H201').Value2 = 'grafische Darstellung d. Organisation, Auflistung der Abteilungen'
$ws.Range('H202').Value2 = 'Liste der Jahresberichte'
$ws.Range('H203').Value2 = 'Berichte Studienbeiträge zu bestimmten Themen'
$ws.Range('H204').Value2 = 'Benutzerordnungen, Datenschutz/-sicherheit,'
$ws.Range('H206').Value2 = 'Login zur Literaturansicht'

# G203's border/alignment formatting was re-applied (xlWrapText + top
# vertical alignment) while typing up the new content next to it.
$g203 = $ws.Range('G203')
$g203.WrapText = $true
$g203.VerticalAlignment = -4160
$g203.Interior.Pattern = -4142

# Scroll position / selection at the moment the sheet was saved.
$excel.ActiveWindow.ScrollRow = 162
$null = $ws.Range('J178').Select()
